$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "C:\Users\320074769\Downloads\Python_Repos\funcSSH\functiondefextractor\test_resource\test_repo\src\"

$ws.Range("A2").Value = $prefix + "CerberusTest.java_testCallMethod"
$ws.Range("A3").Value = $prefix + "CerberusTest.java_testCerebrusWithArguments"
$ws.Range("A4").Value = $prefix + "CerberusTest.java_testCerebruswithOutArguments"
$ws.Range("A5").Value = $prefix + "CerberusTest.java_testCerebruswithWrongArguments"

$cr = [char]13
$lf = [char]10

$ws.Range("B2").Value = "@Test" + $cr + "public void testCallMethod() throws Exception {" + $cr + "assertEquals(Integer.valueOf(0), new Cerberus().call())" + $cr + "}" + $cr + $lf

$ws.Range("B3").Value = "@Test" + $cr + "public void testCerebrusWithArguments() {" + $cr + "getOriginalOutputStream().flush()" + $cr + "Cerberus.main(new String[] { `"CPD`" })" + $cr + "}" + $cr + $lf

$ws.Range("B4").Value = "@Test" + $cr + "public void testCerebruswithOutArguments() {" + $cr + "Cerberus.main(new String[] {})" + $cr + "String expectedOutputString = getCerberusCommandLineUsageString()" + $cr + "assertEquals(expectedOutputString, getModifiedOutputStream().toString())" + $cr + "}" + $cr + $lf

$ws.Range("B5").Value = "@Test" + $cr + "public void testCerebruswithWrongArguments() {" + $cr + "String dummyArgument = `"dummy argument`"" + $cr + "Cerberus.main(new String[] { dummyArgument })" + $cr + "String expectedOutputString = new StringBuilder().append(`"Unmatched argument at index 0: 'dummy argument'`").append(NEW_LINE).append(getCerberusCommandLineUsageString()).toString()" + $cr + "}" + $cr + $lf

# The multi-line values above trigger Excel's automatic row-height recalculation
# (customHeight). The original file has no explicit row heights, so re-run
# AutoFit to drop the explicit height/customHeight markers back to default,
# keeping row formatting untouched (the commit only changed text content).
$ws.UsedRange.EntireRow.AutoFit()
